$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 8324.77
$ws.Range("I33").Value = 87.75
$ws.Range("J33").Value = 21504
$ws.Range("K33").Value = 87.75
$ws.Range("L33").Value = 21504
$ws.Range("M33").Value = 141.25
$ws.Range("N33").Value = -21962

$ws.Range("H70").Value = 1834
$ws.Range("I70").Value = 1502
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 4506
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -4236
$ws.Range("N70").Value = -6540

$ws.Range("H73").Value = 1834
$ws.Range("I73").Value = 1502
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 4506
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -3570
$ws.Range("N73").Value = -7872

$ws.Range("H86").Value = 2131.375
$ws.Range("I86").Value = 1902.65
$ws.Range("K86").Value = 1902.65
$ws.Range("M86").Value = -779.6500000000001

$ws.Range("H89").Value = 2131.375
$ws.Range("I89").Value = 1902.65
$ws.Range("K89").Value = 9513.25
$ws.Range("M89").Value = -3897.25

$ws.Range("H132").Value = 2369.3333
$ws.Range("I132").Value = 1678.9117
$ws.Range("J132").Value = 4503.364
$ws.Range("K132").Value = 5036.7351
$ws.Range("L132").Value = 13510.092
$ws.Range("M132").Value = -2506.7351
$ws.Range("N132").Value = -18570.092

$ws.Range("H137").Value = 3232.8
$ws.Range("I137").Value = 4923.2607
$ws.Range("J137").Value = 945.7059
$ws.Range("K137").Value = 14769.7821
$ws.Range("L137").Value = 2837.1177
$ws.Range("M137").Value = -12219.7821
$ws.Range("N137").Value = -7937.117700000001

$ws.Range("H138").Value = 2612.6875
$ws.Range("I138").Value = 1240.849
$ws.Range("J138").Value = 5305.5557
$ws.Range("K138").Value = 3722.547
$ws.Range("L138").Value = 15916.6671
$ws.Range("M138").Value = 1417.453
$ws.Range("N138").Value = -26196.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 346936
$ws.Range("I32").Value = 2186.3247
$ws.Range("J32").Value = 4771223.5
$ws.Range("K32").Value = 2186.3247
$ws.Range("L32").Value = 4771223.5
$ws.Range("M32").Value = -1899.3247
$ws.Range("N32").Value = -4771797.5

$ws.Range("H63").Value = 4677.6924
$ws.Range("I63").Value = 3060.8
$ws.Range("J63").Value = 5688.25
$ws.Range("K63").Value = 3060.8
$ws.Range("L63").Value = 5688.25
$ws.Range("M63").Value = -2374.8
$ws.Range("N63").Value = -7060.25

$ws.Range("H66").Value = 4677.6924
$ws.Range("I66").Value = 3060.8
$ws.Range("J66").Value = 5688.25
$ws.Range("K66").Value = 15304
$ws.Range("L66").Value = 28441.25
$ws.Range("M66").Value = -11872
$ws.Range("N66").Value = -35305.25

$ws.Range("H74").Value = 358009.03
$ws.Range("I74").Value = 400849.62
$ws.Range("J74").Value = 1004
$ws.Range("K74").Value = 400849.62
$ws.Range("L74").Value = 1004
$ws.Range("M74").Value = -399975.62
$ws.Range("N74").Value = -2752

$ws.Range("H77").Value = 358009.03
$ws.Range("I77").Value = 400849.62
$ws.Range("J77").Value = 1004
$ws.Range("K77").Value = 2004248.1
$ws.Range("L77").Value = 5020
$ws.Range("M77").Value = -1999880.1
$ws.Range("N77").Value = -13756

$ws.Range("H97").Value = 2260.0667
$ws.Range("I97").Value = 1589
$ws.Range("J97").Value = 3602.2
$ws.Range("K97").Value = 1589
$ws.Range("L97").Value = 3602.2
$ws.Range("M97").Value = -1093
$ws.Range("N97").Value = -4594.2

$ws.Range("H122").Value = 25658944
$ws.Range("I122").Value = 29431962
$ws.Range("J122").Value = 2428.4
$ws.Range("K122").Value = 88295886
$ws.Range("L122").Value = 7285.200000000001
$ws.Range("M122").Value = -88293436
$ws.Range("N122").Value = -12185.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9093032
$ws.Range("I31").Value = 11629514
$ws.Range("J31").Value = 3967.75
$ws.Range("K31").Value = 11629514
$ws.Range("L31").Value = 3967.75
$ws.Range("M31").Value = -11629219
$ws.Range("N31").Value = -4557.75

$ws.Range("H34").Value = 9093032
$ws.Range("I34").Value = 11629514
$ws.Range("J34").Value = 3967.75
$ws.Range("K34").Value = 11629514
$ws.Range("L34").Value = 3967.75
$ws.Range("M34").Value = -11629312
$ws.Range("N34").Value = -4371.75

$ws.Range("H132").Value = 39588.117
$ws.Range("I132").Value = 1017.1818
$ws.Range("J132").Value = 251728.25
$ws.Range("K132").Value = 3051.5454
$ws.Range("L132").Value = 755184.75
$ws.Range("M132").Value = -521.5454
$ws.Range("N132").Value = -760244.75

$ws.Range("H134").Value = 19127.418
$ws.Range("I134").Value = 23712.455
$ws.Range("J134").Value = 787.2727
$ws.Range("K134").Value = 71137.365
$ws.Range("L134").Value = 2361.8181
$ws.Range("M134").Value = -68602.365
$ws.Range("N134").Value = -7431.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 889.8
$ws.Range("I122").Value = 825
$ws.Range("J122").Value = 933
$ws.Range("K122").Value = 7425
$ws.Range("L122").Value = 8397
$ws.Range("M122").Value = -4975
$ws.Range("N122").Value = -13297

$ws.Range("H131").Value = 6850203
$ws.Range("J131").Value = 7937380.5
$ws.Range("L131").Value = 23812141.5
$ws.Range("N131").Value = -23822221.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 5675
$ws.Range("I42").Value = 6825
$ws.Range("J42").Value = 5100
$ws.Range("K42").Value = 6825
$ws.Range("L42").Value = 5100
$ws.Range("M42").Value = -6262
$ws.Range("N42").Value = -6226

$ws.Range("H49").Value = 5675
$ws.Range("I49").Value = 6825
$ws.Range("J49").Value = 5100
$ws.Range("K49").Value = 6825
$ws.Range("L49").Value = 5100
$ws.Range("M49").Value = -6678
$ws.Range("N49").Value = -5394

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 17125710
$ws.Range("I132").Value = 24039566
$ws.Range("J132").Value = 5686.905
$ws.Range("K132").Value = 72118698
$ws.Range("L132").Value = 17060.715
$ws.Range("M132").Value = -72116168
$ws.Range("N132").Value = -22120.715
